$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric data in columns B..H for rows 2..13 to the nearest integer,
# matching the "write to disk as integer data" behaviour described in the commit.
for ($r = 2; $r -le 13; $r++) {
    for ($c = 2; $c -le 8; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val) {
            $cell.Value2 = [Math]::Round([double]$val, 0)
        }
    }
}
